# "update brick 1.6 report"
# Disable (Enable=0) the first 12 test cases (case1..case12, rows 2-13) on
# the "data" sheet, and move the view/selection down to the bottom of the
# table (rows 14-25), keeping the existing frozen header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate() | Out-Null

# Flip the "Enable" flag (column B) from 1 to 0 for case1..case12
$ws.Range("B2:B13").Value = 0

# Scroll the frozen (bottom) pane down so row 23 is the first visible row,
# then select B14:B25 (mirrors the author's on-screen selection when they
# saved the file).
$win = $excel.Application.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1

$ws.Range("B14:B25").Select() | Out-Null

# Restore the workbook window position recorded by the author's Excel session.
$excel.Left = -8415
$excel.Top = -16320
